$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.878.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.752.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.58%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.394"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -16.12%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.236.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.783.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.751.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "361.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.568"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.172"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0938"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.58%  "
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "332.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0598"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0259"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.641"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "137.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.53%  "
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.88%  "
